$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The "Gent vs St. Truiden" and "Genk vs Charleroi" matches (rows 30 and
# --- 31) had their match data swapped onto the wrong row. Columns A:E
# --- (Indice, pais, torneio, temporada, data_partida) are unaffected;
# --- rewrite columns F:V on each row with the correct match's data.

# Row 30 should now hold the Genk vs Charleroi data.
$ws.Cells.Item(30, 6).Value = "Genk"
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = "Charleroi"
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 1.47
$ws.Cells.Item(30, 11).Value = "13/08/2023 18:42"
$ws.Cells.Item(30, 12).Value = 1.59
$ws.Cells.Item(30, 13).Value = "20/08/2023 15:53"
$ws.Cells.Item(30, 14).Value = 4.85
$ws.Cells.Item(30, 15).Value = "13/08/2023 18:42"
$ws.Cells.Item(30, 16).Value = 4.57
$ws.Cells.Item(30, 17).Value = "20/08/2023 15:53"
$ws.Cells.Item(30, 18).Value = 5.52
$ws.Cells.Item(30, 19).Value = "13/08/2023 18:42"
$ws.Cells.Item(30, 20).Value = 5.25
$ws.Cells.Item(30, 21).Value = "20/08/2023 15:58"
$ws.Cells.Item(30, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/genk-charleroi/xlfRPeMc/"

# Row 31 should now hold the Gent vs St. Truiden data.
$ws.Cells.Item(31, 6).Value = "Gent"
$ws.Cells.Item(31, 7).Value = 2
$ws.Cells.Item(31, 8).Value = "St. Truiden"
$ws.Cells.Item(31, 9).Value = 2
$ws.Cells.Item(31, 10).Value = 1.4
$ws.Cells.Item(31, 11).Value = "13/08/2023 19:42"
$ws.Cells.Item(31, 12).Value = 1.52
$ws.Cells.Item(31, 13).Value = "20/08/2023 15:57"
$ws.Cells.Item(31, 14).Value = 5
$ws.Cells.Item(31, 15).Value = "13/08/2023 19:42"
$ws.Cells.Item(31, 16).Value = 4.45
$ws.Cells.Item(31, 17).Value = "20/08/2023 15:59"
$ws.Cells.Item(31, 18).Value = 6.31
$ws.Cells.Item(31, 19).Value = "13/08/2023 19:42"
$ws.Cells.Item(31, 20).Value = 6.41
$ws.Cells.Item(31, 21).Value = "20/08/2023 15:59"
$ws.Cells.Item(31, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/gent-st-truiden/SxmvNg6G/"

# --- Append a new match row (row 122) at the end of the sheet: Westerlo vs
# --- Anderlecht, scraped 01-12-2023 20:45.
$newRow = 122

# Clone formatting from the row above (index column + date column carry
# special styles) before filling in the new row's values.
$ws.Cells.Item($newRow - 1, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($newRow - 1, 5).Copy() | Out-Null
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = 121
$ws.Cells.Item($newRow, 2).Value = "belgium"
$ws.Cells.Item($newRow, 3).Value = "jupiler-pro-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45261.86458333334
$ws.Cells.Item($newRow, 6).Value = "Westerlo"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Anderlecht"
$ws.Cells.Item($newRow, 9).Value = 3
$ws.Cells.Item($newRow, 10).Value = 4.06
$ws.Cells.Item($newRow, 11).Value = "26/11/2023 18:43"
$ws.Cells.Item($newRow, 12).Value = 5.14
$ws.Cells.Item($newRow, 13).Value = "01/12/2023 20:44"
$ws.Cells.Item($newRow, 14).Value = 4.05
$ws.Cells.Item($newRow, 15).Value = "26/11/2023 18:43"
$ws.Cells.Item($newRow, 16).Value = 3.83
$ws.Cells.Item($newRow, 17).Value = "01/12/2023 20:42"
$ws.Cells.Item($newRow, 18).Value = 1.75
$ws.Cells.Item($newRow, 19).Value = "26/11/2023 18:43"
$ws.Cells.Item($newRow, 20).Value = 1.72
$ws.Cells.Item($newRow, 21).Value = "01/12/2023 20:35"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/westerlo-anderlecht/E7hwOl90/"
